# BackupCodes.xlsx -- "Added CCDI 5 scripts and C3DC structure"
#
# The "sec-codes" sheet holds a flat list of one-time backup codes in
# column A. This refreshes the first batch of codes (rows 2-4) with newly
# issued values, drops the now-unused codes that used to sit in rows
# 10-12, and leaves the remaining codes (rows 13-15) untouched in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sec-codes")

# Rotate in the new codes for rows 2-4.
$ws.Range("A2").Value = "C9DEXVFAR31A"
$ws.Range("A3").Value = "05ANDJ337D9B"
$ws.Range("A4").Value = "SAS5DZQK4GHR"

# Remove the stale codes that used to live in rows 10-12 without
# shifting the codes below them (rows 13-15 keep their row numbers).
$ws.Range("A10:A12").ClearContents()

# Matches the author's last on-screen selection after the edit.
$ws.Range("A6").Select()
